$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 3170.2
$ws.Range("I113").Value = 3414.818
$ws.Range("J113").Value = 2497.5
$ws.Range("K113").Value = 3414.818
$ws.Range("L113").Value = 2497.5
$ws.Range("M113").Value = -160.8180000000002
$ws.Range("N113").Value = -9005.5

# Row 137
$ws.Range("H137").Value = 40572.117
$ws.Range("I137").Value = 1560
$ws.Range("K137").Value = 4680
$ws.Range("M137").Value = -2130


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 11632912
$ws.Range("I2").Value = 11632912
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 11632912
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -11632799

# Row 23
$ws.Range("H23").Value = 75006.5
$ws.Range("I23").Value = 70006
$ws.Range("J23").Value = 80007
$ws.Range("K23").Value = 70006
$ws.Range("L23").Value = 80007
$ws.Range("M23").Value = -69747
$ws.Range("N23").Value = -80525

# Row 61
$ws.Range("H61").Value = 3687.4
$ws.Range("I61").Value = 2772.2222
$ws.Range("J61").Value = 4436.1816
$ws.Range("K61").Value = 2772.2222
$ws.Range("L61").Value = 4436.1816
$ws.Range("M61").Value = -2560.2222
$ws.Range("N61").Value = -4860.1816

# Row 74
$ws.Range("H74").Value = 2785.7144
$ws.Range("I74").Value = 850
$ws.Range("K74").Value = 850
$ws.Range("M74").Value = 24

# Row 77
$ws.Range("H77").Value = 2785.7144
$ws.Range("I77").Value = 850
$ws.Range("K77").Value = 4250
$ws.Range("M77").Value = 118

# Row 110
$ws.Range("H110").Value = 897.8333
$ws.Range("I110").Value = 897.4
$ws.Range("J110").Value = 900
$ws.Range("K110").Value = 897.4
$ws.Range("L110").Value = 900
$ws.Range("M110").Value = 1147.6
$ws.Range("N110").Value = -4990

# Row 116
$ws.Range("H116").Value = 11632912
$ws.Range("I116").Value = 11632912
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 11632912
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -11630618

# Row 136
$ws.Range("H136").Value = 3687.4
$ws.Range("I136").Value = 2772.2222
$ws.Range("J136").Value = 4436.1816
$ws.Range("K136").Value = 8316.6666
$ws.Range("L136").Value = 13308.5448
$ws.Range("M136").Value = -5766.6666
$ws.Range("N136").Value = -18408.5448


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 11632912
$ws.Range("I3").Value = 11632912
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 11632912
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -11632798

# Row 86
$ws.Range("H86").Value = 1001500
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877

# Row 89
$ws.Range("H89").Value = 1001500
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 15000
$ws.Range("M89").Value = -9384

# Row 132
$ws.Range("H132").Value = 70733.336
$ws.Range("J132").Value = 70733.336
$ws.Range("L132").Value = 70733.336
$ws.Range("N132").Value = -80853.336

# Row 134
$ws.Range("H134").Value = 3095.7778
$ws.Range("I134").Value = 2851.7144
$ws.Range("K134").Value = 8555.143199999999
$ws.Range("M134").Value = -6020.143199999999


$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 10417326
$ws.Range("I22").Value = 594
$ws.Range("J22").Value = 20834058
$ws.Range("K22").Value = 594
$ws.Range("L22").Value = 20834058
$ws.Range("M22").Value = -244
$ws.Range("N22").Value = -20834758

# Row 122
$ws.Range("H122").Value = 1031.0714
$ws.Range("J122").Value = 800
$ws.Range("L122").Value = 2400
$ws.Range("N122").Value = -7300

# Row 134
$ws.Range("H134").Value = 1929.6364
$ws.Range("I134").Value = 935.5
$ws.Range("K134").Value = 2806.5
$ws.Range("M134").Value = -271.5


$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 7864129
$ws.Range("I11").Value = 6823806
$ws.Range("K11").Value = 6823806
$ws.Range("M11").Value = -6823667

# Row 122
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 126
$ws.Range("H126").Value = 3706972
$ws.Range("I126").Value = 7940109.5
$ws.Range("J126").Value = 2976.75
$ws.Range("K126").Value = 23820328.5
$ws.Range("L126").Value = 8930.25
$ws.Range("M126").Value = -23817858.5
$ws.Range("N126").Value = -13870.25

# Row 132
$ws.Range("H132").Value = 1071014.4
$ws.Range("I132").Value = 1674025.2
$ws.Range("J132").Value = 4148.923
$ws.Range("K132").Value = 5022075.6
$ws.Range("L132").Value = 12446.769
$ws.Range("M132").Value = -5019545.6
$ws.Range("N132").Value = -17506.769


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4156.8423
$ws.Range("I16").Value = 4268.3335
$ws.Range("J16").Value = 2150
$ws.Range("K16").Value = 4268.3335
$ws.Range("L16").Value = 2150
$ws.Range("M16").Value = -4098.3335
$ws.Range("N16").Value = -2490

# Row 40
$ws.Range("H40").Value = 9923.5
$ws.Range("I40").Value = 4847.75
$ws.Range("K40").Value = 4847.75
$ws.Range("M40").Value = -4711.75

# Row 46
$ws.Range("H46").Value = 888.7778
$ws.Range("I46").Value = 540
$ws.Range("J46").Value = 1324.75
$ws.Range("K46").Value = 540
$ws.Range("L46").Value = 1324.75
$ws.Range("M46").Value = -352
$ws.Range("N46").Value = -1700.75

# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0

# Row 122
$ws.Range("H122").Value = 13571.571
$ws.Range("I122").Value = 13751
$ws.Range("K122").Value = 41253
$ws.Range("M122").Value = -38803

# Row 136
$ws.Range("H136").Value = 5646.067
$ws.Range("I136").Value = 3332.8333
$ws.Range("K136").Value = 9998.499899999999
$ws.Range("M136").Value = -7448.499899999999


$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

# Row 65
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

# Row 107
$ws.Range("H107").Value = 948.4
$ws.Range("I107").Value = 380.875
$ws.Range("K107").Value = 1142.625
$ws.Range("M107").Value = 777.375

# Row 122
$ws.Range("H122").Value = 50505.19
$ws.Range("I122").Value = 65236.8
$ws.Range("J122").Value = 1399.8334
$ws.Range("K122").Value = 195710.4
$ws.Range("L122").Value = 4199.5002
$ws.Range("M122").Value = -193260.4
$ws.Range("N122").Value = -9099.5002

# Row 136
$ws.Range("H136").Value = 27781146
$ws.Range("I136").Value = 46299372
$ws.Range("K136").Value = 138898116
$ws.Range("M136").Value = -138895566

